# Apply the diff: reshuffled "Canto Coral IV" syllabus sheet.
# Final layout: rows 10-23 get corrected content; row 24 (old) is removed
# (dimension shrinks from A1:C24 to A1:C23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$elisabeth = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$objEnLong = '1. Approximate student of their vocal apparatus, at the level of spoken and sung words, in their individual and collective expression (choral). 2. Encourage the student to experience the choral repertoire and its role in the development of musical language. 3. Provide the student with the opportunity to vocal music together, with technical learning parameters such as tuning, precision, balance, phrasing etc.'
$semestral = 'Semestral'
$resumptionShortEn = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing.Vocal classification. Breathing into the singing. Placing the emission in Bocca Chiusa. Choral singing in unison. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.'
$date = '01/01/2017'
$resumptionLongEn = 'Resumption and improvement of repertoires and technical assumptions of previous discipline Choir Singing. Vocal classification.Breathing into the singing - Exercises for locating low and average breathing. Support and air column. Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. Choral singing in canon. Choral singing in other polyphonic formations. Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission.'
$metodo = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas – a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
$criterio = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
$naoTem = 'Não tem'
$requisito = "8800011 -  Canto Coral III  (Requisito)`n"

# Row 10: Objetivos: -> Elisabeth text (was long PT objectives paragraph)
$ws.Range("B10").Value = $elisabeth
$ws.Range("C10").Value = $elisabeth

# Row 11: Objectives: -> unchanged (long EN objectives paragraph)
$ws.Range("B11").Value = $objEnLong
$ws.Range("C11").Value = $objEnLong

# Row 12: Docentes responsáveis: label only, no content - unchanged

# Row 13: used to hold the (unlabeled) Elisabeth text; now becomes
# "Programa resumido:" / "Semestral", with a 60pt custom row height.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = $semestral
$ws.Range("C13").Value = $semestral
$ws.Rows(13).RowHeight = 60

# Row 14: Short syllabus -> unchanged content, but label moves here
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = $resumptionShortEn
$ws.Range("C14").Value = $resumptionShortEn

# Row 15: Programa: -> date string "01/01/2017", height grows to 120
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = $date
$ws.Range("C15").Value = $date
$ws.Rows(15).RowHeight = 120

# Row 16: Syllabus: -> unchanged long EN syllabus text
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = $resumptionLongEn
$ws.Range("C16").Value = $resumptionLongEn

# Row 17: becomes "Avaliação:" label only (content cleared, height reset)
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows(17).AutoFit()

# Row 18: becomes "Método:" / Elisabeth text, 60pt height
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = $elisabeth
$ws.Range("C18").Value = $elisabeth
$ws.Rows(18).RowHeight = 60

# Row 19: becomes "Critério:" / the old "Método" content
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20: becomes "Norma de recuperação:" / the old "Critério" content
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21: becomes "Bibliografia:" / "Não tem", height reset to default
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = $naoTem
$ws.Range("C21").Value = $naoTem
$ws.Rows(21).AutoFit()

# Row 22: becomes "Requisitos:" label only (content cleared, height reset)
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").Clear()
$ws.Rows(22).AutoFit()

# Row 23 previously held "Requisitos:" alone; row 24 held the requisite
# text. Move that text up into row 23 and delete the now-superfluous row 24.
$ws.Range("A23").Clear()
$ws.Range("B23").Value = $requisito
$ws.Range("C23").Value = $requisito
$ws.Rows(23).RowHeight = 30

$ws.Rows(24).Delete()
